$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------
# Sheet "展览" (Exhibition) - sheet1
# ----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Simple numeric F-column updates (ticket/attendance counts)
$ws1.Range("F8").Value = 2035   # was 2030
$ws1.Range("F10").Value = 39   # was 38
$ws1.Range("F12").Value = 1604   # was 1602
$ws1.Range("F13").Value = 1604   # was 1602
$ws1.Range("F14").Value = 1333   # was 1332
$ws1.Range("F16").Value = 1391   # was 1390
$ws1.Range("F19").Value = 11   # was 10
$ws1.Range("F20").Value = 489   # was 482
$ws1.Range("F23").Value = 7080   # was 7075
$ws1.Range("F24").Value = 7689   # was 7685
$ws1.Range("F28").Value = 492   # was 491
$ws1.Range("F29").Value = 85   # was 84
$ws1.Range("F41").Value = 699   # was 698

# Row 36-39 content shift: row 36 (cancelled event) is removed and
# replaced by what followed it; a new row lands at 39.
$ws1.Range("B36").NumberFormat = "@"
$ws1.Range("B36").Value = '2024-10-04'
$ws1.Range("C36").Value = '北京·第五人格only同人展'
$ws1.Range("D36").Value = '北花园路1号 超级蜂巢'
$ws1.Range("E36").Value = '2024.10.04 10:00-10.04 17:00'
$ws1.Range("F36").Value = 1394
$ws1.Range("G36").Value = 68
$ws1.Range("H36").Value = 'https://show.bilibili.com/platform/detail.html?id=89309'
$ws1.Range("I36").Value = '//i0.hdslb.com/bfs/openplatform/202407/4XsICpa71721046044404.jpeg'

$ws1.Range("B37").NumberFormat = "@"
$ws1.Range("B37").Value = '2024-10-05'
$ws1.Range("C37").Value = '北京·咒术回战同人Only2.0'
$ws1.Range("D37").Value = '安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园'
$ws1.Range("E37").Value = '2024.10.05 09:30-10.05 17:00'
$ws1.Range("F37").Value = 21
$ws1.Range("G37").Value = 65
$ws1.Range("H37").Value = 'https://show.bilibili.com/platform/detail.html?id=91628'
$ws1.Range("I37").Value = '//i0.hdslb.com/bfs/openplatform/202408/IsJo7aU61724405528082.jpeg'

$ws1.Range("B38").NumberFormat = "@"
$ws1.Range("B38").Value = '2024-10-05'
$ws1.Range("C38").Value = '北京·马娘ONLY2'
$ws1.Range("D38").Value = '永外高庄138号 北京大红门国际会展中心'
$ws1.Range("E38").Value = '2024.10.05 10:00-10.05 17:00'
$ws1.Range("F38").Value = 215
$ws1.Range("G38").Value = 75
$ws1.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=89334'
$ws1.Range("I38").Value = '//i1.hdslb.com/bfs/openplatform/202408/30C9r9Qz1724639124911.png'

$ws1.Range("B39").NumberFormat = "@"
$ws1.Range("B39").Value = '2024-10-06'
$ws1.Range("C39").Value = '北京·Hi Fun 全忍界秋季运动会 火影同人ONLY x 北投购物公园潮街 '
$ws1.Range("D39").Value = '安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园'
$ws1.Range("E39").Value = '2024.10.06 11:00-10.06 18:00'
$ws1.Range("F39").Value = 1
$ws1.Range("G39").Value = 68
$ws1.Range("H39").Value = 'https://show.bilibili.com/platform/detail.html?id=91930'
$ws1.Range("I39").Value = '//i2.hdslb.com/bfs/openplatform/202409/f7nTqmEI1725439502652.jpeg'

# ----------------------------------------------------------------
# Sheet "演出" (Performance) - sheet2
# ----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16   # was 15

# ----------------------------------------------------------------
# Sheet "本地生活" (Local life) - sheet3
# ----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 178   # was 177
$ws3.Range("F3").Value = 2586   # was 2584
$ws3.Range("F5").Value = 127   # was 125
$ws3.Range("F6").Value = 4   # was 3

# ----------------------------------------------------------------
# Sheet "全部类型" (All types) - sheet4
# ----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 178   # was 177
$ws4.Range("F4").Value = 16   # was 15
$ws4.Range("F5").Value = 127   # was 125
$ws4.Range("F10").Value = 2035   # was 2030
$ws4.Range("F11").Value = 39   # was 38
$ws4.Range("F13").Value = 1604   # was 1602
$ws4.Range("F14").Value = 1604   # was 1602
$ws4.Range("F15").Value = 4   # was 3
$ws4.Range("F16").Value = 1333   # was 1332
$ws4.Range("F18").Value = 11   # was 10
$ws4.Range("F19").Value = 489   # was 482
$ws4.Range("F23").Value = 7080   # was 7075
$ws4.Range("F24").Value = 7689   # was 7685
$ws4.Range("F28").Value = 85   # was 84
$ws4.Range("F34").Value = 215   # was 216
$ws4.Range("F39").Value = 699   # was 698
